# ---------------------------------------------------------------------------
# chore: adapt column header formatting to respective input file names (#7)
#
# The workbook is an "AHB diff" export: columns A-J describe the old
# ("_old") format version, L-U the new ("_new") format version (K holds the
# literal "diff" marker). This change renames those generic "_old"/"_new"
# header suffixes to the concrete format-version tags they represent
# (FV2310 / FV2404), freezes the header row, and wraps the used range in a
# real Excel Table ("Table1") whose column names mirror the header row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header row (row 1): "<Name>_old" -> "<Name>_FV2310",
#    "<Name>_new" -> "<Name>_FV2404". The "diff" column header (K1) is left
#    untouched since it carries neither suffix.
$headerMap = @{
    "A1" = "Segmentname_FV2310"
    "B1" = "Segmentgruppe_FV2310"
    "C1" = "Segment_FV2310"
    "D1" = "Datenelement_FV2310"
    "E1" = "Segment ID_FV2310"
    "F1" = "Code_FV2310"
    "G1" = "Qualifier_FV2310"
    "H1" = "Beschreibung_FV2310"
    "I1" = "Bedingungsausdruck_FV2310"
    "J1" = "Bedingung_FV2310"
    "L1" = "Segmentname_FV2404"
    "M1" = "Segmentgruppe_FV2404"
    "N1" = "Segment_FV2404"
    "O1" = "Datenelement_FV2404"
    "P1" = "Segment ID_FV2404"
    "Q1" = "Code_FV2404"
    "R1" = "Qualifier_FV2404"
    "S1" = "Beschreibung_FV2404"
    "T1" = "Bedingungsausdruck_FV2404"
    "U1" = "Bedingung_FV2404"
}

foreach ($addr in $headerMap.Keys) {
    $ws.Range($addr).Value = $headerMap[$addr]
}

# 2) Freeze the header row so it stays visible while scrolling: select the
#    first cell below the header and freeze panes above/left of it.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3) Turn the whole used range (including the just-renamed header row) into
#    an Excel Table named "Table1" with a header row and an AutoFilter.
$rng = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add(1, $rng, 0, 1)
$tbl.Name = "Table1"
